$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (pushes existing rows 25-58 down to 26-59)
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the new data record
$ws.Cells.Item(25, 1).Value2  = 11
$ws.Cells.Item(25, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value2  = "Bíobío"
$ws.Cells.Item(25, 4).Value2  = 44979
$ws.Cells.Item(25, 5).Value2  = 8
$ws.Cells.Item(25, 6).Value2  = 100112030
$ws.Cells.Item(25, 7).Value2  = "Poroto granado"
$ws.Cells.Item(25, 8).Value2  = "Sin especificar"
$ws.Cells.Item(25, 9).Value2  = "Primera"
$ws.Cells.Item(25, 10).Value2 = 100
$ws.Cells.Item(25, 11).Value2 = 25000
$ws.Cells.Item(25, 12).Value2 = 26000
$ws.Cells.Item(25, 13).Value2 = 25500
$ws.Cells.Item(25, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value2 = 1020
$ws.Cells.Item(25, 17).Value2 = 25
$ws.Cells.Item(25, 18).Value2 = "Hortaliza"

# Match the date cell formatting used throughout column D (numFmt for dates)
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
